$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 2
$ws.Range("B2").Value = "w"
$ws.Range("C2").Value = "@gmail.com"
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1
$ws.Range("H2").Value = 71.86

# Row 3
$ws.Range("A3").Value = 3
$ws.Range("B3").Value = "w"
$ws.Range("C3").Value = "@gmail.com"
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 16.09
